# Update TPM-derived NATMI ligand-receptor metrics (Thbs1-Tnfrsf11b) for rows 2-19, columns G:T.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 34.70905366666667
$ws.Range("H2").Value = 104.127161
$ws.Range("I2").Value = 0.01618617796956752
$ws.Range("J2").Value = 0.01618617796956752
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2262196666666666
$ws.Range("N2").Value = 0.6786589999999999
$ws.Range("O2").Value = 0.03145179203784564
$ws.Range("P2").Value = 0.03145179203784564
$ws.Range("Q2").Value = 7.851870550788777
$ws.Range("R2").Value = 70.66683495709898
$ws.Range("S2").Value = 0.0005090843033863961
$ws.Range("T2").Value = 0.0005090843033863961

$ws.Range("G3").Value = 34.70905366666667
$ws.Range("H3").Value = 104.127161
$ws.Range("I3").Value = 0.01618617796956752
$ws.Range("J3").Value = 0.01618617796956752
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.93109
$ws.Range("N3").Value = 20.79327
$ws.Range("O3").Value = 0.9636438974901603
$ws.Range("P3").Value = 0.9636438974901604
$ws.Range("Q3").Value = 240.5715747784967
$ws.Range("R3").Value = 2165.14417300647
$ws.Range("S3").Value = 0.01559771162406341
$ws.Range("T3").Value = 0.01559771162406341

$ws.Range("G4").Value = 34.70905366666667
$ws.Range("H4").Value = 104.127161
$ws.Range("I4").Value = 0.01618617796956752
$ws.Range("J4").Value = 0.01618617796956752
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03527466666666667
$ws.Range("N4").Value = 0.105824
$ws.Range("O4").Value = 0.004904310471994002
$ws.Range("P4").Value = 0.004904310471994003
$ws.Range("Q4").Value = 1.224350298407111
$ws.Range("R4").Value = 11.019152685664
$ws.Range("S4").Value = 7.938204211770859 / 100000
$ws.Range("T4").Value = 7.938204211770861 / 100000

$ws.Range("G5").Value = 95.50314333333334
$ws.Range("H5").Value = 286.50943
$ws.Range("I5").Value = 0.04453681997475516
$ws.Range("J5").Value = 0.04453681997475516
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2262196666666666
$ws.Range("N5").Value = 0.6786589999999999
$ws.Range("O5").Value = 0.03145179203784564
$ws.Range("P5").Value = 0.03145179203784564
$ws.Range("Q5").Value = 21.60468925048555
$ws.Range("R5").Value = 194.44220325437
$ws.Range("S5").Value = 0.001400762799872969
$ws.Range("T5").Value = 0.001400762799872969

$ws.Range("G6").Value = 95.50314333333334
$ws.Range("H6").Value = 286.50943
$ws.Range("I6").Value = 0.04453681997475516
$ws.Range("J6").Value = 0.04453681997475516
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.93109
$ws.Range("N6").Value = 20.79327
$ws.Range("O6").Value = 0.9636438974901603
$ws.Range("P6").Value = 0.9636438974901604
$ws.Range("Q6").Value = 661.9408817262334
$ws.Range("R6").Value = 5957.4679355361
$ws.Range("S6").Value = 0.04291763478229069
$ws.Range("T6").Value = 0.04291763478229069

$ws.Range("G7").Value = 95.50314333333334
$ws.Range("H7").Value = 286.50943
$ws.Range("I7").Value = 0.04453681997475516
$ws.Range("J7").Value = 0.04453681997475516
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03527466666666667
$ws.Range("N7").Value = 0.105824
$ws.Range("O7").Value = 0.004904310471994002
$ws.Range("P7").Value = 0.004904310471994003
$ws.Range("Q7").Value = 3.368841546702223
$ws.Range("R7").Value = 30.31957392032
$ws.Range("S7").Value = 0.0002184223925915034
$ws.Range("T7").Value = 0.0002184223925915034

$ws.Range("G8").Value = 1028.132161333333
$ws.Range("H8").Value = 3084.396484
$ws.Range("I8").Value = 0.4794579045397415
$ws.Range("J8").Value = 0.4794579045397416
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2262196666666666
$ws.Range("N8").Value = 0.6786589999999999
$ws.Range("O8").Value = 0.03145179203784564
$ws.Range("P8").Value = 0.03145179203784564
$ws.Range("Q8").Value = 232.5837148261062
$ws.Range("R8").Value = 2093.253433434956
$ws.Range("S8").Value = 0.0150798103044852
$ws.Range("T8").Value = 0.0150798103044852

$ws.Range("G9").Value = 1028.132161333333
$ws.Range("H9").Value = 3084.396484
$ws.Range("I9").Value = 0.4794579045397415
$ws.Range("J9").Value = 0.4794579045397416
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.93109
$ws.Range("N9").Value = 20.79327
$ws.Range("O9").Value = 0.9636438974901603
$ws.Range("P9").Value = 0.9636438974901604
$ws.Range("Q9").Value = 7126.076542095852
$ws.Range("R9").Value = 64134.68887886267
$ws.Range("S9").Value = 0.4620266838131417
$ws.Range("T9").Value = 0.4620266838131419

$ws.Range("G10").Value = 1028.132161333333
$ws.Range("H10").Value = 3084.396484
$ws.Range("I10").Value = 0.4794579045397415
$ws.Range("J10").Value = 0.4794579045397416
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03527466666666667
$ws.Range("N10").Value = 0.105824
$ws.Range("O10").Value = 0.004904310471994002
$ws.Range("P10").Value = 0.004904310471994003
$ws.Range("Q10").Value = 36.26701928031289
$ws.Range("R10").Value = 326.403173522816
$ws.Range("S10").Value = 0.002351410422114555
$ws.Range("T10").Value = 0.002351410422114556

$ws.Range("G11").Value = 47.840114
$ws.Range("H11").Value = 143.520342
$ws.Range("I11").Value = 0.02230970071166346
$ws.Range("J11").Value = 0.02230970071166346
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.2262196666666666
$ws.Range("N11").Value = 0.6786589999999999
$ws.Range("O11").Value = 0.03145179203784564
$ws.Range("P11").Value = 0.03145179203784564
$ws.Range("Q11").Value = 10.82237464237533
$ws.Range("R11").Value = 97.40137178137799
$ws.Range("S11").Value = 0.0007016800672098159
$ws.Range("T11").Value = 0.0007016800672098159

$ws.Range("G12").Value = 47.840114
$ws.Range("H12").Value = 143.520342
$ws.Range("I12").Value = 0.02230970071166346
$ws.Range("J12").Value = 0.02230970071166346
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 6.93109
$ws.Range("N12").Value = 20.79327
$ws.Range("O12").Value = 0.9636438974901603
$ws.Range("P12").Value = 0.9636438974901604
$ws.Range("Q12").Value = 331.58413574426
$ws.Range("R12").Value = 2984.25722169834
$ws.Range("S12").Value = 0.02149860694562638
$ws.Range("T12").Value = 0.02149860694562638

$ws.Range("G13").Value = 47.840114
$ws.Range("H13").Value = 143.520342
$ws.Range("I13").Value = 0.02230970071166346
$ws.Range("J13").Value = 0.02230970071166346
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.03527466666666667
$ws.Range("N13").Value = 0.105824
$ws.Range("O13").Value = 0.004904310471994002
$ws.Range("P13").Value = 0.004904310471994003
$ws.Range("Q13").Value = 1.687544074645333
$ws.Range("R13").Value = 15.187896671808
$ws.Range("S13").Value = 0.0001094136988272631
$ws.Range("T13").Value = 0.0001094136988272632

$ws.Range("G14").Value = 428.234253
$ws.Range("H14").Value = 1284.702759
$ws.Range("I14").Value = 0.19970224190755
$ws.Range("J14").Value = 0.19970224190755
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.2262196666666666
$ws.Range("N14").Value = 0.6786589999999999
$ws.Range("O14").Value = 0.03145179203784564
$ws.Range("P14").Value = 0.03145179203784564
$ws.Range("Q14").Value = 96.87500996890898
$ws.Range("R14").Value = 871.8750897201809
$ws.Range("S14").Value = 0.006280993381967805
$ws.Range("T14").Value = 0.006280993381967805

$ws.Range("G15").Value = 428.234253
$ws.Range("H15").Value = 1284.702759
$ws.Range("I15").Value = 0.19970224190755
$ws.Range("J15").Value = 0.19970224190755
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 6.93109
$ws.Range("N15").Value = 20.79327
$ws.Range("O15").Value = 0.9636438974901603
$ws.Range("P15").Value = 0.9636438974901604
$ws.Range("Q15").Value = 2968.13014862577
$ws.Range("R15").Value = 26713.17133763193
$ws.Range("S15").Value = 0.1924418467293143
$ws.Range("T15").Value = 0.1924418467293143

$ws.Range("G16").Value = 428.234253
$ws.Range("H16").Value = 1284.702759
$ws.Range("I16").Value = 0.19970224190755
$ws.Range("J16").Value = 0.19970224190755
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03527466666666667
$ws.Range("N16").Value = 0.105824
$ws.Range("O16").Value = 0.004904310471994002
$ws.Range("P16").Value = 0.004904310471994003
$ws.Range("Q16").Value = 15.105820529824
$ws.Range("R16").Value = 135.952384768416
$ws.Range("S16").Value = 0.0009794017962678769
$ws.Range("T16").Value = 0.0009794017962678771

$ws.Range("G17").Value = 509.945048
$ws.Range("H17").Value = 1529.835144
$ws.Range("I17").Value = 0.2378071548967224
$ws.Range("J17").Value = 0.2378071548967224
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.2262196666666666
$ws.Range("N17").Value = 0.6786589999999999
$ws.Range("O17").Value = 0.03145179203784564
$ws.Range("P17").Value = 0.03145179203784564
$ws.Range("Q17").Value = 115.3595987768773
$ws.Range("R17").Value = 1038.236388991896
$ws.Range("S17").Value = 0.007479461180923457
$ws.Range("T17").Value = 0.007479461180923458

$ws.Range("G18").Value = 509.945048
$ws.Range("H18").Value = 1529.835144
$ws.Range("I18").Value = 0.2378071548967224
$ws.Range("J18").Value = 0.2378071548967224
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 6.93109
$ws.Range("N18").Value = 20.79327
$ws.Range("O18").Value = 0.9636438974901603
$ws.Range("P18").Value = 0.9636438974901604
$ws.Range("Q18").Value = 3534.475022742321
$ws.Range("R18").Value = 31810.27520468088
$ws.Range("S18").Value = 0.2291614135957238
$ws.Range("T18").Value = 0.2291614135957238

$ws.Range("G19").Value = 509.945048
$ws.Range("H19").Value = 1529.835144
$ws.Range("I19").Value = 0.2378071548967224
$ws.Range("J19").Value = 0.2378071548967224
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.03527466666666667
$ws.Range("N19").Value = 0.105824
$ws.Range("O19").Value = 0.004904310471994002
$ws.Range("P19").Value = 0.004904310471994003
$ws.Range("Q19").Value = 17.98814158651734
$ws.Range("R19").Value = 161.893274278656
$ws.Range("S19").Value = 0.001166280120075095
$ws.Range("T19").Value = 0.001166280120075095

